# Add a new trial row ("Pa_Fri_43dpf_GroupC_n5_20200711_1550") right after the
# existing "SF_Sat_14dpf_GroupA_n5_20200613_1205" / secondary_factor row (row 9),
# pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 10; everything from old row 10 downward
# (including the two blank-but-styled rows near the bottom) shifts down by one.
$ws.Rows.Item(10).Insert() | Out-Null

# Copy the formatting of the row that is now 11 (the old row 10) onto the new,
# still-blank row 10 so fills/fonts/column widths line up with its neighbours.
$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("A10:D10").PasteSpecial(-4122) | Out-Null

# Fill in the new trial's data.
$ws.Cells.Item(10, 1).Value = "Pa_Fri_43dpf_GroupC_n5_20200711_1550"
$ws.Cells.Item(10, 2).Value = "bkgSub_options.secondary_factor"
$ws.Cells.Item(10, 3).Value = 3
$ws.Cells.Item(10, 4).Value = ""

# Keep the cursor/view roughly where the original author left it.
$ws.Range("D28").Select() | Out-Null
